# Auto-generated Excel COM-interop script to update cryptos.xlsx data
# Applies the per-cell changes described in the commit diff (price/volume refresh,
# plus three row swaps: Toncoin<->InjectiveProtocol, WEMIXToken<->ARBITRUM, Algorand->RocketPoolETH).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.796.43"
$ws.Range("E2").Value = "  -1.00%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.592.50"
$ws.Range("E3").Value = "  -0.45%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.86"
$ws.Range("E5").Value = "  -0.72%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.41"
$ws.Range("E6").Value = "  -1.86%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.591"
$ws.Range("E7").Value = "  -1.29%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.575"
$ws.Range("E9").Value = "  -1.44%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.50"
$ws.Range("E10").Value = "  -0.40%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.02"
$ws.Range("E11").Value = "  -1.24%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0838"
$ws.Range("E12").Value = "  -0.38%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.05"
$ws.Range("E13").Value = "  -3.83%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.988.80"
$ws.Range("E14").Value = "  -0.39%  "

# Row 15
$ws.Range("E15").Value = "  +1.21%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.597.80"
$ws.Range("E16").Value = "  -0.08%  "

# Row 17
$ws.Range("E17").Value = "  +0.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.71"
$ws.Range("E18").Value = "  -1.01%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "45.906.27"
$ws.Range("E19").Value = "  -1.19%  "

# Row 20
$ws.Range("E20").Value = "  -0.49%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.68"
$ws.Range("E21").Value = "  -0.36%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.51"
$ws.Range("E22").Value = "  -4.40%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "287.50"
$ws.Range("E23").Value = "  +12.66%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.89"
$ws.Range("E24").Value = "  +2.34%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.02"
$ws.Range("E25").Value = "  -1.98%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  +0.58%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.11"
$ws.Range("E27").Value = "  +2.32%  "

# Row 28
$ws.Range("E28").Value = "  +0.08%  "

# Row 29
$ws.Range("E29").Value = "  +0.38%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.64"
$ws.Range("E30").Value = "  +1.11%  "

# Row 31
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.35"
$ws.Range("E31").Value = "  -3.59%  "

# Row 32
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.18"
$ws.Range("E32").Value = "  -4.27%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.31"
$ws.Range("E33").Value = "  +2.96%  "

# Row 34
$ws.Range("E34").Value = "  -1.79%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.25"
$ws.Range("E35").Value = "  +3.53%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.21"
$ws.Range("E36").Value = "  -3.31%  "

# Row 37
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.81"
$ws.Range("E37").Value = "  -2.98%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0832"
$ws.Range("E38").Value = "  +0.06%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.122"
$ws.Range("E39").Value = "  +4.23%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.122"
$ws.Range("E40").Value = "  +0.17%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.47"
$ws.Range("E41").Value = "  -4.69%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0325"
$ws.Range("E42").Value = "  +0.40%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.99"
$ws.Range("E44").Value = "  -4.48%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.99"
$ws.Range("E45").Value = "  -3.66%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.119.80"
$ws.Range("E46").Value = "  +3.85%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("E47").Value = "  -0.06%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.33"
$ws.Range("E48").Value = "  +3.36%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.22"
$ws.Range("E49").Value = "  -0.73%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "108.10"
$ws.Range("E50").Value = "  -1.40%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.842.55"
$ws.Range("E51").Value = "  -0.71%  "
